$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40. This shifts the existing rows 40-148
# (and all their values/styles) down to rows 41-149, leaving a blank
# (but style-inheriting) row 40 ready to be populated.
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with a new data record (same shape as the
# surrounding rows), with the date and volume values called out by the diff.
$ws.Range("A40").Value = 3
$ws.Range("B40").Value = "Femacal de La Calera"
$ws.Range("C40").Value = "Coquimbo"
$ws.Range("D40").Value = 44414
$ws.Range("E40").Value = 5
$ws.Range("F40").Value = 100112039
$ws.Range("G40").Value = "Ciboulette"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 160
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = 1500
$ws.Range("N40").Value = "$/docena de atados"
$ws.Range("O40").Value = "Provincia de Quillota"
$ws.Range("P40").Value = 500
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = "Hortaliza"
